# Adds the new "달바글로벌" (Dalba Global) sheet with its date/remn_amt
# history, and appends the latest (2025-10-28 / serial 45958) data point
# as a new row to each of the five pre-existing ticker sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append row 101 to the five existing sheets (same date, new amount).
# ---------------------------------------------------------------------
$existingAdds = "485883;366242;185804;257613;533239" -split ";"

for ($i = 0; $i -lt 5; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $lastRow = $ws.Cells.Item(100, 1).Row
    $newRow = $lastRow + 1

    # Carry the column-A date style down from the row above so the new
    # cell keeps the same number format (style index 2 in the original
    # workbook) instead of minting a near-duplicate style.
    $ws.Cells.Item($lastRow, 1).Copy($ws.Cells.Item($newRow, 1))

    $ws.Cells.Item($newRow, 1).Value = 45958
    $ws.Cells.Item($newRow, 2).Value = [int]$existingAdds[$i]
}

# ---------------------------------------------------------------------
# 2. Add the brand-new "달바글로벌" sheet as the last tab.
# ---------------------------------------------------------------------
$firstWs = $wb.Worksheets.Item(1)
$lastWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$dalba = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastWs)
$dalba.Name = "달바글로벌"

# Copy the header cells (A1:B1) from an existing sheet so the new sheet
# picks up the same bold/border/centered header style instead of a
# freshly minted one.
$firstWs.Range("A1:B1").Copy($dalba.Range("A1:B1"))
$dalba.Cells.Item(1, 1).Value = "date"
$dalba.Cells.Item(1, 2).Value = "remn_amt"

# Copy the column-A date style (style index 2) down across all 100 data
# rows before filling in values, same reasoning as step 1.
$firstWs.Cells.Item(2, 1).Copy($dalba.Range("A2:A101"))

$dalbaCsv = "45807,251;45810,2687;45812,4312;45813,4274;45817,4423;45818,6294;45819,6849;45820,8390;45821,11198;45824,11673;45825,11867;45826,13234;45827,13740;45828,13942;45831,15023;45832,15753;45833,12586;45834,12266;45835,12137;45838,12933;45839,12679;45840,12817;45841,12982;45842,12815;45845,14176;45846,16153;45847,17178;45848,17680;45849,18674;45852,19546;45853,19939;45854,19319;45855,20124;45856,18787;45859,19135;45860,18775;45861,18934;45862,18507;45863,18686;45866,20374;45867,20586;45868,20054;45869,20004;45870,19163;45873,19951;45874,21075;45875,21680;45876,25748;45877,25163;45880,29028;45881,27531;45882,27409;45883,28523;45887,27782;45888,28715;45889,28557;45890,37533;45891,45210;45894,49704;45895,49251;45896,49954;45897,52945;45898,54247;45901,56297;45902,56293;45903,55189;45904,56188;45905,58256;45908,58711;45909,61242;45910,61364;45911,61489;45912,60110;45915,61288;45916,62955;45917,62569;45918,64068;45919,60363;45922,59891;45923,59581;45924,58441;45925,57951;45926,56105;45929,57447;45930,57161;45931,57097;45932,58205;45940,56425;45943,56208;45944,55401;45945,56394;45946,54403;45947,57423;45950,56545;45951,62741;45952,63817;45953,65132;45954,62679;45957,60948;45958,58330"

$dalbaRows = $dalbaCsv -split ";"
$n = $dalbaRows.Length

$data = New-Object 'object[,]' $n, 2
for ($i = 0; $i -lt $n; $i++) {
    $parts = $dalbaRows[$i] -split ","
    $data[$i, 0] = [int]$parts[0]
    $data[$i, 1] = [int]$parts[1]
}

$dalba.Range("A2:B101").Value = $data
